# "Open Account Test Added"
#
# The workbook already has an "OpenAccount" sheet (2 columns: Customer /
# Currency, one data row). This adds a third column ("alertText") and
# replaces the sample data row with the new Open-Account test fixture
# (customer "Harry Potter", currency "Dollar", and the expected success
# message), mirroring the existing "AddCustomer" sheet's alertText column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpenAccount")
$ws.Activate()

# New header cell for the third column.
$ws.Range("C1").Value = "alertText"

# Replace the sample row with the Open-Account fixture. C2 is written
# first so the new shared strings are interned in the same order as the
# target workbook (alertText value, then customer name, then currency).
$ws.Range("C2").Value = "Account created successfully"
$ws.Range("A2").Value = "Harry Potter"
$ws.Range("B2").Value = "Dollar"

# Widen the name / alertText columns to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 12.625
$ws.Columns.Item(3).ColumnWidth = 27.140625

# Leave the selection where the author left it when saving.
$ws.Range("B3").Select()
